$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 317; this shifts existing rows 317-370 down to 318-371,
# carrying their full A:T contents with them automatically.
$ws.Rows.Item(317).Insert()

# The newly inserted row 317 is blank; populate the constant columns that are shared by
# every record in this subset (Mercado ID, Mercado, Region, Codreg, Tipo, Producto ID,
# Producto, Categoria ID, Categoria, Origen).
$ws.Range("A317").Value = 7
$ws.Range("B317").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C317").Value = 'Ñuble'
$ws.Range("E317").Value = 16
$ws.Range("F317").Value = 'Fruta'
$ws.Range("G317").Value = 100108
$ws.Range("H317").Value = 'Tropicales y subtropicales'
$ws.Range("I317").Value = 100108005
$ws.Range("J317").Value = 'Piña'
$ws.Range("R317").Value = 'Ecuador'

# Row 316: Fecha, Variedad, Calidad, Volumen, Precio minimo/maximo/promedio, Unidad, Precio $/Kg, Kg/unidad
$ws.Range("D316").Value = 45180
$ws.Range("K316").Value = 'Caramelo'
$ws.Range("L316").Value = 'Segunda'
$ws.Range("M316").Value = 40
$ws.Range("N316").Value = 23000
$ws.Range("O316").Value = 23000
$ws.Range("P316").Value = 23000
$ws.Range("Q316").Value = '$/caja 14 unidades'
$ws.Range("S316").Value = 1643
$ws.Range("T316").Value = 14

# Row 317: Fecha, Variedad, Calidad, Volumen, Precio minimo/maximo/promedio, Unidad, Precio $/Kg, Kg/unidad
$ws.Range("D317").Value = 44340
$ws.Range("K317").Value = 'Caramelo'
$ws.Range("L317").Value = 'Segunda'
$ws.Range("M317").Value = 120
$ws.Range("N317").Value = 15500
$ws.Range("O317").Value = 16000
$ws.Range("P317").Value = 15750
$ws.Range("Q317").Value = '$/caja 14 unidades'
$ws.Range("S317").Value = 1125
$ws.Range("T317").Value = 14

# Row 318: Fecha, Variedad, Calidad, Volumen, Precio minimo/maximo/promedio, Unidad, Precio $/Kg, Kg/unidad
$ws.Range("D318").Value = 44175
$ws.Range("K318").Value = 'Caramelo'
$ws.Range("L318").Value = 'Primera'
$ws.Range("M318").Value = 55
$ws.Range("N318").Value = 19000
$ws.Range("O318").Value = 20000
$ws.Range("P318").Value = 19545
$ws.Range("Q318").Value = '$/caja 12 unidades'
$ws.Range("S318").Value = 1629
$ws.Range("T318").Value = 12

# Row 319: Fecha, Variedad, Calidad, Volumen, Precio minimo/maximo/promedio, Unidad, Precio $/Kg, Kg/unidad
$ws.Range("D319").Value = 44175
$ws.Range("K319").Value = 'Caramelo'
$ws.Range("L319").Value = 'Segunda'
$ws.Range("M319").Value = 40
$ws.Range("N319").Value = 19500
$ws.Range("O319").Value = 20000
$ws.Range("P319").Value = 19750
$ws.Range("Q319").Value = '$/caja 14 unidades'
$ws.Range("S319").Value = 1411
$ws.Range("T319").Value = 14

# Row 320: Fecha, Variedad, Calidad, Volumen, Precio minimo/maximo/promedio, Unidad, Precio $/Kg, Kg/unidad
$ws.Range("D320").Value = 44649
$ws.Range("K320").Value = 'Caramelo'
$ws.Range("L320").Value = 'Segunda'
$ws.Range("M320").Value = 100
$ws.Range("N320").Value = 17000
$ws.Range("O320").Value = 18000
$ws.Range("P320").Value = 17500
$ws.Range("Q320").Value = '$/caja 14 unidades'
$ws.Range("S320").Value = 1250
$ws.Range("T320").Value = 14

# Row 321: Fecha, Variedad, Calidad, Volumen, Precio minimo/maximo/promedio, Unidad, Precio $/Kg, Kg/unidad
$ws.Range("D321").Value = 44273
$ws.Range("K321").Value = 'Caramelo'
$ws.Range("L321").Value = 'Segunda'
$ws.Range("M321").Value = 120
$ws.Range("N321").Value = 15000
$ws.Range("O321").Value = 16000
$ws.Range("P321").Value = 15500
$ws.Range("Q321").Value = '$/caja 14 unidades'
$ws.Range("S321").Value = 1107
$ws.Range("T321").Value = 14

# Row 322: Fecha, Variedad, Calidad, Volumen, Precio minimo/maximo/promedio, Unidad, Precio $/Kg, Kg/unidad
$ws.Range("D322").Value = 44799
$ws.Range("K322").Value = 'Caramelo'
$ws.Range("L322").Value = 'Segunda'
$ws.Range("M322").Value = 120
$ws.Range("N322").Value = 19000
$ws.Range("O322").Value = 20000
$ws.Range("P322").Value = 19500
$ws.Range("Q322").Value = '$/caja 14 unidades'
$ws.Range("S322").Value = 1393
$ws.Range("T322").Value = 14

# Row 323: Fecha, Variedad, Calidad, Volumen, Precio minimo/maximo/promedio, Unidad, Precio $/Kg, Kg/unidad
$ws.Range("D323").Value = 44321
$ws.Range("K323").Value = 'Caramelo'
$ws.Range("L323").Value = 'Segunda'
$ws.Range("M323").Value = 120
$ws.Range("N323").Value = 15000
$ws.Range("O323").Value = 16000
$ws.Range("P323").Value = 15500
$ws.Range("Q323").Value = '$/caja 14 unidades'
$ws.Range("S323").Value = 1107
$ws.Range("T323").Value = 14

# Row 324: Fecha, Variedad, Calidad, Volumen, Precio minimo/maximo/promedio, Unidad, Precio $/Kg, Kg/unidad
$ws.Range("D324").Value = 45093
$ws.Range("K324").Value = 'Caramelo'
$ws.Range("L324").Value = 'Primera'
$ws.Range("M324").Value = 60
$ws.Range("N324").Value = 22000
$ws.Range("O324").Value = 22000
$ws.Range("P324").Value = 22000
$ws.Range("Q324").Value = '$/caja 12 unidades'
$ws.Range("S324").Value = 1833
$ws.Range("T324").Value = 12

# Row 325: Fecha, Variedad, Calidad, Volumen, Precio minimo/maximo/promedio, Unidad, Precio $/Kg, Kg/unidad
$ws.Range("D325").Value = 44291
$ws.Range("K325").Value = 'Caramelo'
$ws.Range("L325").Value = 'Segunda'
$ws.Range("M325").Value = 120
$ws.Range("N325").Value = 15000
$ws.Range("O325").Value = 16000
$ws.Range("P325").Value = 15500
$ws.Range("Q325").Value = '$/caja 14 unidades'
$ws.Range("S325").Value = 1107
$ws.Range("T325").Value = 14

# Row 326: Fecha, Variedad, Calidad, Volumen, Precio minimo/maximo/promedio, Unidad, Precio $/Kg, Kg/unidad
$ws.Range("D326").Value = 44414
$ws.Range("K326").Value = 'Caramelo'
$ws.Range("L326").Value = 'Primera'
$ws.Range("M326").Value = 100
$ws.Range("N326").Value = 18000
$ws.Range("O326").Value = 19000
$ws.Range("P326").Value = 18500
$ws.Range("Q326").Value = '$/caja 12 unidades'
$ws.Range("S326").Value = 1542
$ws.Range("T326").Value = 12

# Row 327: Fecha, Variedad, Calidad, Volumen, Precio minimo/maximo/promedio, Unidad, Precio $/Kg, Kg/unidad
$ws.Range("D327").Value = 44414
$ws.Range("K327").Value = 'Caramelo'
$ws.Range("L327").Value = 'Segunda'
$ws.Range("M327").Value = 100
$ws.Range("N327").Value = 18000
$ws.Range("O327").Value = 19000
$ws.Range("P327").Value = 18500
$ws.Range("Q327").Value = '$/caja 14 unidades'
$ws.Range("S327").Value = 1321
$ws.Range("T327").Value = 14

# Row 328: Fecha, Variedad, Calidad, Volumen, Precio minimo/maximo/promedio, Unidad, Precio $/Kg, Kg/unidad
$ws.Range("D328").Value = 44414
$ws.Range("K328").Value = 'Caramelo'
$ws.Range("L328").Value = 'Tercera'
$ws.Range("M328").Value = 40
$ws.Range("N328").Value = 18000
$ws.Range("O328").Value = 18000
$ws.Range("P328").Value = 18000
$ws.Range("Q328").Value = '$/caja 16 unidades'
$ws.Range("S328").Value = 1125
$ws.Range("T328").Value = 16

# Row 329: Fecha, Variedad, Calidad, Volumen, Precio minimo/maximo/promedio, Unidad, Precio $/Kg, Kg/unidad
$ws.Range("D329").Value = 44181
$ws.Range("K329").Value = 'Caramelo'
$ws.Range("L329").Value = 'Primera'
$ws.Range("M329").Value = 45
$ws.Range("N329").Value = 16500
$ws.Range("O329").Value = 17000
$ws.Range("P329").Value = 16778
$ws.Range("Q329").Value = '$/caja 12 unidades'
$ws.Range("S329").Value = 1398
$ws.Range("T329").Value = 12

# Row 330: Fecha, Variedad, Calidad, Volumen, Precio minimo/maximo/promedio, Unidad, Precio $/Kg, Kg/unidad
$ws.Range("D330").Value = 44181
$ws.Range("K330").Value = 'Caramelo'
$ws.Range("L330").Value = 'Segunda'
$ws.Range("M330").Value = 60
$ws.Range("N330").Value = 16500
$ws.Range("O330").Value = 17000
$ws.Range("P330").Value = 16708
$ws.Range("Q330").Value = '$/caja 14 unidades'
$ws.Range("S330").Value = 1193
$ws.Range("T330").Value = 14

# Row 331: Fecha, Variedad, Calidad, Volumen, Precio minimo/maximo/promedio, Unidad, Precio $/Kg, Kg/unidad
$ws.Range("D331").Value = 45075
$ws.Range("K331").Value = 'Caramelo'
$ws.Range("L331").Value = 'Primera'
$ws.Range("M331").Value = 30
$ws.Range("N331").Value = 15000
$ws.Range("O331").Value = 15000
$ws.Range("P331").Value = 15000
$ws.Range("Q331").Value = '$/caja 12 unidades'
$ws.Range("S331").Value = 1250
$ws.Range("T331").Value = 12

# Row 332: Fecha, Variedad, Calidad, Volumen, Precio minimo/maximo/promedio, Unidad, Precio $/Kg, Kg/unidad
$ws.Range("D332").Value = 44960
$ws.Range("K332").Value = 'Caramelo'
$ws.Range("L332").Value = 'Segunda'
$ws.Range("M332").Value = 30
$ws.Range("N332").Value = 18000
$ws.Range("O332").Value = 18000
$ws.Range("P332").Value = 18000
$ws.Range("Q332").Value = '$/caja 14 unidades'
$ws.Range("S332").Value = 1286
$ws.Range("T332").Value = 14

# Row 333: Fecha, Variedad, Calidad, Volumen, Precio minimo/maximo/promedio, Unidad, Precio $/Kg, Kg/unidad
$ws.Range("D333").Value = 44883
$ws.Range("K333").Value = 'Caramelo'
$ws.Range("L333").Value = 'Segunda'
$ws.Range("M333").Value = 100
$ws.Range("N333").Value = 27000
$ws.Range("O333").Value = 28000
$ws.Range("P333").Value = 27500
$ws.Range("Q333").Value = '$/caja 14 unidades'
$ws.Range("S333").Value = 1964
$ws.Range("T333").Value = 14

# Row 334: Fecha, Variedad, Calidad, Volumen, Precio minimo/maximo/promedio, Unidad, Precio $/Kg, Kg/unidad
$ws.Range("D334").Value = 44309
$ws.Range("K334").Value = 'Caramelo'
$ws.Range("L334").Value = 'Segunda'
$ws.Range("M334").Value = 60
$ws.Range("N334").Value = 15000
$ws.Range("O334").Value = 16000
$ws.Range("P334").Value = 15500
$ws.Range("Q334").Value = '$/caja 14 unidades'
$ws.Range("S334").Value = 1107
$ws.Range("T334").Value = 14

# Row 335: Fecha, Variedad, Calidad, Volumen, Precio minimo/maximo/promedio, Unidad, Precio $/Kg, Kg/unidad
$ws.Range("D335").Value = 44672
$ws.Range("K335").Value = 'Caramelo'
$ws.Range("L335").Value = 'Segunda'
$ws.Range("M335").Value = 120
$ws.Range("N335").Value = 14000
$ws.Range("O335").Value = 15000
$ws.Range("P335").Value = 14500
$ws.Range("Q335").Value = '$/caja 14 unidades'
$ws.Range("S335").Value = 1036
$ws.Range("T335").Value = 14

# Row 336: Fecha, Variedad, Calidad, Volumen, Precio minimo/maximo/promedio, Unidad, Precio $/Kg, Kg/unidad
$ws.Range("D336").Value = 44650
$ws.Range("K336").Value = 'Caramelo'
$ws.Range("L336").Value = 'Segunda'
$ws.Range("M336").Value = 120
$ws.Range("N336").Value = 17000
$ws.Range("O336").Value = 18000
$ws.Range("P336").Value = 17500
$ws.Range("Q336").Value = '$/caja 14 unidades'
$ws.Range("S336").Value = 1250
$ws.Range("T336").Value = 14

# Row 337: Fecha, Variedad, Calidad, Volumen, Precio minimo/maximo/promedio, Unidad, Precio $/Kg, Kg/unidad
$ws.Range("D337").Value = 44699
$ws.Range("K337").Value = 'Caramelo'
$ws.Range("L337").Value = 'Segunda'
$ws.Range("M337").Value = 120
$ws.Range("N337").Value = 17000
$ws.Range("O337").Value = 18000
$ws.Range("P337").Value = 17500
$ws.Range("Q337").Value = '$/caja 14 unidades'
$ws.Range("S337").Value = 1250
$ws.Range("T337").Value = 14

# Row 338: Fecha, Variedad, Calidad, Volumen, Precio minimo/maximo/promedio, Unidad, Precio $/Kg, Kg/unidad
$ws.Range("D338").Value = 44413
$ws.Range("K338").Value = 'Caramelo'
$ws.Range("L338").Value = 'Primera'
$ws.Range("M338").Value = 60
$ws.Range("N338").Value = 18000
$ws.Range("O338").Value = 19000
$ws.Range("P338").Value = 18500
$ws.Range("Q338").Value = '$/caja 12 unidades'
$ws.Range("S338").Value = 1542
$ws.Range("T338").Value = 12

# Row 339: Fecha, Variedad, Calidad, Volumen, Precio minimo/maximo/promedio, Unidad, Precio $/Kg, Kg/unidad
$ws.Range("D339").Value = 44413
$ws.Range("K339").Value = 'Caramelo'
$ws.Range("L339").Value = 'Segunda'
$ws.Range("M339").Value = 60
$ws.Range("N339").Value = 18000
$ws.Range("O339").Value = 19000
$ws.Range("P339").Value = 18500
$ws.Range("Q339").Value = '$/caja 14 unidades'
$ws.Range("S339").Value = 1321
$ws.Range("T339").Value = 14

# Row 340: Fecha, Variedad, Calidad, Volumen, Precio minimo/maximo/promedio, Unidad, Precio $/Kg, Kg/unidad
$ws.Range("D340").Value = 44453
$ws.Range("K340").Value = 'Caramelo'
$ws.Range("L340").Value = 'Primera'
$ws.Range("M340").Value = 60
$ws.Range("N340").Value = 19000
$ws.Range("O340").Value = 20000
$ws.Range("P340").Value = 19500
$ws.Range("Q340").Value = '$/caja 12 unidades'
$ws.Range("S340").Value = 1625
$ws.Range("T340").Value = 12

# Row 341: Fecha, Variedad, Calidad, Volumen, Precio minimo/maximo/promedio, Unidad, Precio $/Kg, Kg/unidad
$ws.Range("D341").Value = 44453
$ws.Range("K341").Value = 'Caramelo'
$ws.Range("L341").Value = 'Segunda'
$ws.Range("M341").Value = 60
$ws.Range("N341").Value = 19000
$ws.Range("O341").Value = 20000
$ws.Range("P341").Value = 19500
$ws.Range("Q341").Value = '$/caja 14 unidades'
$ws.Range("S341").Value = 1393
$ws.Range("T341").Value = 14

# Row 342: Fecha, Variedad, Calidad, Volumen, Precio minimo/maximo/promedio, Unidad, Precio $/Kg, Kg/unidad
$ws.Range("D342").Value = 44319
$ws.Range("K342").Value = 'Caramelo'
$ws.Range("L342").Value = 'Segunda'
$ws.Range("M342").Value = 60
$ws.Range("N342").Value = 16000
$ws.Range("O342").Value = 17000
$ws.Range("P342").Value = 16500
$ws.Range("Q342").Value = '$/caja 14 unidades'
$ws.Range("S342").Value = 1179
$ws.Range("T342").Value = 14

# Row 343: Fecha, Variedad, Calidad, Volumen, Precio minimo/maximo/promedio, Unidad, Precio $/Kg, Kg/unidad
$ws.Range("D343").Value = 45166
$ws.Range("K343").Value = 'Caramelo'
$ws.Range("L343").Value = 'Segunda'
$ws.Range("M343").Value = 60
$ws.Range("N343").Value = 22000
$ws.Range("O343").Value = 22000
$ws.Range("P343").Value = 22000
$ws.Range("Q343").Value = '$/caja 14 unidades'
$ws.Range("S343").Value = 1571
$ws.Range("T343").Value = 14

# Row 344: Fecha, Variedad, Calidad, Volumen, Precio minimo/maximo/promedio, Unidad, Precio $/Kg, Kg/unidad
$ws.Range("D344").Value = 44336
$ws.Range("K344").Value = 'Caramelo'
$ws.Range("L344").Value = 'Segunda'
$ws.Range("M344").Value = 120
$ws.Range("N344").Value = 15500
$ws.Range("O344").Value = 16000
$ws.Range("P344").Value = 15750
$ws.Range("Q344").Value = '$/caja 14 unidades'
$ws.Range("S344").Value = 1125
$ws.Range("T344").Value = 14

# Row 345: Fecha, Variedad, Calidad, Volumen, Precio minimo/maximo/promedio, Unidad, Precio $/Kg, Kg/unidad
$ws.Range("D345").Value = 44659
$ws.Range("K345").Value = 'Caramelo'
$ws.Range("L345").Value = 'Segunda'
$ws.Range("M345").Value = 120
$ws.Range("N345").Value = 15000
$ws.Range("O345").Value = 16000
$ws.Range("P345").Value = 15500
$ws.Range("Q345").Value = '$/caja 14 unidades'
$ws.Range("S345").Value = 1107
$ws.Range("T345").Value = 14

# Row 346: Fecha, Variedad, Calidad, Volumen, Precio minimo/maximo/promedio, Unidad, Precio $/Kg, Kg/unidad
$ws.Range("D346").Value = 44924
$ws.Range("K346").Value = 'Caramelo'
$ws.Range("L346").Value = 'Segunda'
$ws.Range("M346").Value = 120
$ws.Range("N346").Value = 19000
$ws.Range("O346").Value = 20000
$ws.Range("P346").Value = 19500
$ws.Range("Q346").Value = '$/caja 14 unidades'
$ws.Range("S346").Value = 1393
$ws.Range("T346").Value = 14

# Row 347: Fecha, Variedad, Calidad, Volumen, Precio minimo/maximo/promedio, Unidad, Precio $/Kg, Kg/unidad
$ws.Range("D347").Value = 44880
$ws.Range("K347").Value = 'Caramelo'
$ws.Range("L347").Value = 'Segunda'
$ws.Range("M347").Value = 120
$ws.Range("N347").Value = 29000
$ws.Range("O347").Value = 30000
$ws.Range("P347").Value = 29500
$ws.Range("Q347").Value = '$/caja 14 unidades'
$ws.Range("S347").Value = 2107
$ws.Range("T347").Value = 14

# Row 348: Fecha, Variedad, Calidad, Volumen, Precio minimo/maximo/promedio, Unidad, Precio $/Kg, Kg/unidad
$ws.Range("D348").Value = 45117
$ws.Range("K348").Value = 'Caramelo'
$ws.Range("L348").Value = 'Primera'
$ws.Range("M348").Value = 80
$ws.Range("N348").Value = 24000
$ws.Range("O348").Value = 24000
$ws.Range("P348").Value = 24000
$ws.Range("Q348").Value = '$/caja 12 unidades'
$ws.Range("S348").Value = 2000
$ws.Range("T348").Value = 12

# Row 349: Fecha, Variedad, Calidad, Volumen, Precio minimo/maximo/promedio, Unidad, Precio $/Kg, Kg/unidad
$ws.Range("D349").Value = 44637
$ws.Range("K349").Value = 'Caramelo'
$ws.Range("L349").Value = 'Segunda'
$ws.Range("M349").Value = 60
$ws.Range("N349").Value = 17000
$ws.Range("O349").Value = 18000
$ws.Range("P349").Value = 17500
$ws.Range("Q349").Value = '$/caja 14 unidades'
$ws.Range("S349").Value = 1250
$ws.Range("T349").Value = 14

# Row 350: Fecha, Variedad, Calidad, Volumen, Precio minimo/maximo/promedio, Unidad, Precio $/Kg, Kg/unidad
$ws.Range("D350").Value = 44362
$ws.Range("K350").Value = 'Caramelo'
$ws.Range("L350").Value = 'Segunda'
$ws.Range("M350").Value = 120
$ws.Range("N350").Value = 15000
$ws.Range("O350").Value = 16000
$ws.Range("P350").Value = 15500
$ws.Range("Q350").Value = '$/caja 14 unidades'
$ws.Range("S350").Value = 1107
$ws.Range("T350").Value = 14

# Row 351: Fecha, Variedad, Calidad, Volumen, Precio minimo/maximo/promedio, Unidad, Precio $/Kg, Kg/unidad
$ws.Range("D351").Value = 45063
$ws.Range("K351").Value = 'Caramelo'
$ws.Range("L351").Value = 'Primera'
$ws.Range("M351").Value = 50
$ws.Range("N351").Value = 15000
$ws.Range("O351").Value = 15000
$ws.Range("P351").Value = 15000
$ws.Range("Q351").Value = '$/caja 12 unidades'
$ws.Range("S351").Value = 1250
$ws.Range("T351").Value = 12

# Row 352: Fecha, Variedad, Calidad, Volumen, Precio minimo/maximo/promedio, Unidad, Precio $/Kg, Kg/unidad
$ws.Range("D352").Value = 45063
$ws.Range("K352").Value = 'Caramelo'
$ws.Range("L352").Value = 'Segunda'
$ws.Range("M352").Value = 50
$ws.Range("N352").Value = 14000
$ws.Range("O352").Value = 14000
$ws.Range("P352").Value = 14000
$ws.Range("Q352").Value = '$/caja 14 unidades'
$ws.Range("S352").Value = 1000
$ws.Range("T352").Value = 14

# Row 353: Fecha, Variedad, Calidad, Volumen, Precio minimo/maximo/promedio, Unidad, Precio $/Kg, Kg/unidad
$ws.Range("D353").Value = 44792
$ws.Range("K353").Value = 'Caramelo'
$ws.Range("L353").Value = 'Segunda'
$ws.Range("M353").Value = 100
$ws.Range("N353").Value = 19000
$ws.Range("O353").Value = 20000
$ws.Range("P353").Value = 19500
$ws.Range("Q353").Value = '$/caja 14 unidades'
$ws.Range("S353").Value = 1393
$ws.Range("T353").Value = 14

# Row 354: Fecha, Variedad, Calidad, Volumen, Precio minimo/maximo/promedio, Unidad, Precio $/Kg, Kg/unidad
$ws.Range("D354").Value = 44557
$ws.Range("K354").Value = 'Caramelo'
$ws.Range("L354").Value = 'Segunda'
$ws.Range("M354").Value = 120
$ws.Range("N354").Value = 17000
$ws.Range("O354").Value = 18000
$ws.Range("P354").Value = 17500
$ws.Range("Q354").Value = '$/caja 14 unidades'
$ws.Range("S354").Value = 1250
$ws.Range("T354").Value = 14

# Row 355: Fecha, Variedad, Calidad, Volumen, Precio minimo/maximo/promedio, Unidad, Precio $/Kg, Kg/unidad
$ws.Range("D355").Value = 44848
$ws.Range("K355").Value = 'Caramelo'
$ws.Range("L355").Value = 'Segunda'
$ws.Range("M355").Value = 100
$ws.Range("N355").Value = 21000
$ws.Range("O355").Value = 22000
$ws.Range("P355").Value = 21500
$ws.Range("Q355").Value = '$/caja 14 unidades'
$ws.Range("S355").Value = 1536
$ws.Range("T355").Value = 14

# Row 356: Fecha, Variedad, Calidad, Volumen, Precio minimo/maximo/promedio, Unidad, Precio $/Kg, Kg/unidad
$ws.Range("D356").Value = 44635
$ws.Range("K356").Value = 'Caramelo'
$ws.Range("L356").Value = 'Segunda'
$ws.Range("M356").Value = 120
$ws.Range("N356").Value = 17000
$ws.Range("O356").Value = 18000
$ws.Range("P356").Value = 17500
$ws.Range("Q356").Value = '$/caja 14 unidades'
$ws.Range("S356").Value = 1250
$ws.Range("T356").Value = 14

# Row 357: Fecha, Variedad, Calidad, Volumen, Precio minimo/maximo/promedio, Unidad, Precio $/Kg, Kg/unidad
$ws.Range("D357").Value = 44385
$ws.Range("K357").Value = 'Caramelo'
$ws.Range("L357").Value = 'Segunda'
$ws.Range("M357").Value = 120
$ws.Range("N357").Value = 17000
$ws.Range("O357").Value = 18000
$ws.Range("P357").Value = 17500
$ws.Range("Q357").Value = '$/caja 14 unidades'
$ws.Range("S357").Value = 1250
$ws.Range("T357").Value = 14

# Row 358: Fecha, Variedad, Calidad, Volumen, Precio minimo/maximo/promedio, Unidad, Precio $/Kg, Kg/unidad
$ws.Range("D358").Value = 44678
$ws.Range("K358").Value = 'Caramelo'
$ws.Range("L358").Value = 'Segunda'
$ws.Range("M358").Value = 120
$ws.Range("N358").Value = 16000
$ws.Range("O358").Value = 17000
$ws.Range("P358").Value = 16500
$ws.Range("Q358").Value = '$/caja 14 unidades'
$ws.Range("S358").Value = 1179
$ws.Range("T358").Value = 14

# Row 359: Fecha, Variedad, Calidad, Volumen, Precio minimo/maximo/promedio, Unidad, Precio $/Kg, Kg/unidad
$ws.Range("D359").Value = 44194
$ws.Range("K359").Value = 'Caramelo'
$ws.Range("L359").Value = 'Segunda'
$ws.Range("M359").Value = 120
$ws.Range("N359").Value = 14000
$ws.Range("O359").Value = 15000
$ws.Range("P359").Value = 14500
$ws.Range("Q359").Value = '$/caja 14 unidades'
$ws.Range("S359").Value = 1036
$ws.Range("T359").Value = 14

# Row 360: Fecha, Variedad, Calidad, Volumen, Precio minimo/maximo/promedio, Unidad, Precio $/Kg, Kg/unidad
$ws.Range("D360").Value = 44771
$ws.Range("K360").Value = 'Caramelo'
$ws.Range("L360").Value = 'Segunda'
$ws.Range("M360").Value = 120
$ws.Range("N360").Value = 20000
$ws.Range("O360").Value = 21000
$ws.Range("P360").Value = 20500
$ws.Range("Q360").Value = '$/caja 14 unidades'
$ws.Range("S360").Value = 1464
$ws.Range("T360").Value = 14

# Row 361: Fecha, Variedad, Calidad, Volumen, Precio minimo/maximo/promedio, Unidad, Precio $/Kg, Kg/unidad
$ws.Range("D361").Value = 44784
$ws.Range("K361").Value = 'Caramelo'
$ws.Range("L361").Value = 'Segunda'
$ws.Range("M361").Value = 120
$ws.Range("N361").Value = 19000
$ws.Range("O361").Value = 20000
$ws.Range("P361").Value = 19500
$ws.Range("Q361").Value = '$/caja 14 unidades'
$ws.Range("S361").Value = 1393
$ws.Range("T361").Value = 14

# Row 362: Fecha, Variedad, Calidad, Volumen, Precio minimo/maximo/promedio, Unidad, Precio $/Kg, Kg/unidad
$ws.Range("D362").Value = 44813
$ws.Range("K362").Value = 'Caramelo'
$ws.Range("L362").Value = 'Segunda'
$ws.Range("M362").Value = 120
$ws.Range("N362").Value = 20000
$ws.Range("O362").Value = 21000
$ws.Range("P362").Value = 20500
$ws.Range("Q362").Value = '$/caja 14 unidades'
$ws.Range("S362").Value = 1464
$ws.Range("T362").Value = 14

# Row 363: Fecha, Variedad, Calidad, Volumen, Precio minimo/maximo/promedio, Unidad, Precio $/Kg, Kg/unidad
$ws.Range("D363").Value = 44638
$ws.Range("K363").Value = 'Caramelo'
$ws.Range("L363").Value = 'Segunda'
$ws.Range("M363").Value = 60
$ws.Range("N363").Value = 17000
$ws.Range("O363").Value = 18000
$ws.Range("P363").Value = 17500
$ws.Range("Q363").Value = '$/caja 14 unidades'
$ws.Range("S363").Value = 1250
$ws.Range("T363").Value = 14

# Row 364: Fecha, Variedad, Calidad, Volumen, Precio minimo/maximo/promedio, Unidad, Precio $/Kg, Kg/unidad
$ws.Range("D364").Value = 44271
$ws.Range("K364").Value = 'Caramelo'
$ws.Range("L364").Value = 'Segunda'
$ws.Range("M364").Value = 60
$ws.Range("N364").Value = 15500
$ws.Range("O364").Value = 16000
$ws.Range("P364").Value = 15750
$ws.Range("Q364").Value = '$/caja 14 unidades'
$ws.Range("S364").Value = 1125
$ws.Range("T364").Value = 14

# Row 365: Fecha, Variedad, Calidad, Volumen, Precio minimo/maximo/promedio, Unidad, Precio $/Kg, Kg/unidad
$ws.Range("D365").Value = 44251
$ws.Range("K365").Value = 'Caramelo'
$ws.Range("L365").Value = 'Primera'
$ws.Range("M365").Value = 100
$ws.Range("N365").Value = 14000
$ws.Range("O365").Value = 15000
$ws.Range("P365").Value = 14650
$ws.Range("Q365").Value = '$/caja 12 unidades'
$ws.Range("S365").Value = 1221
$ws.Range("T365").Value = 12

# Row 366: Fecha, Variedad, Calidad, Volumen, Precio minimo/maximo/promedio, Unidad, Precio $/Kg, Kg/unidad
$ws.Range("D366").Value = 44286
$ws.Range("K366").Value = 'Caramelo'
$ws.Range("L366").Value = 'Segunda'
$ws.Range("M366").Value = 120
$ws.Range("N366").Value = 15500
$ws.Range("O366").Value = 16000
$ws.Range("P366").Value = 15750
$ws.Range("Q366").Value = '$/caja 14 unidades'
$ws.Range("S366").Value = 1125
$ws.Range("T366").Value = 14

# Row 367: Fecha, Variedad, Calidad, Volumen, Precio minimo/maximo/promedio, Unidad, Precio $/Kg, Kg/unidad
$ws.Range("D367").Value = 44343
$ws.Range("K367").Value = 'Caramelo'
$ws.Range("L367").Value = 'Segunda'
$ws.Range("M367").Value = 120
$ws.Range("N367").Value = 15500
$ws.Range("O367").Value = 16000
$ws.Range("P367").Value = 15750
$ws.Range("Q367").Value = '$/caja 14 unidades'
$ws.Range("S367").Value = 1125
$ws.Range("T367").Value = 14

# Row 368: Fecha, Variedad, Calidad, Volumen, Precio minimo/maximo/promedio, Unidad, Precio $/Kg, Kg/unidad
$ws.Range("D368").Value = 44754
$ws.Range("K368").Value = 'Sin especificar'
$ws.Range("L368").Value = 'Segunda'
$ws.Range("M368").Value = 120
$ws.Range("N368").Value = 18000
$ws.Range("O368").Value = 19000
$ws.Range("P368").Value = 18500
$ws.Range("Q368").Value = '$/caja 14 unidades'
$ws.Range("S368").Value = 1321
$ws.Range("T368").Value = 14

# Row 369: Fecha, Variedad, Calidad, Volumen, Precio minimo/maximo/promedio, Unidad, Precio $/Kg, Kg/unidad
$ws.Range("D369").Value = 44490
$ws.Range("K369").Value = 'Caramelo'
$ws.Range("L369").Value = 'Segunda'
$ws.Range("M369").Value = 60
$ws.Range("N369").Value = 20000
$ws.Range("O369").Value = 21000
$ws.Range("P369").Value = 20500
$ws.Range("Q369").Value = '$/caja 14 unidades'
$ws.Range("S369").Value = 1464
$ws.Range("T369").Value = 14

# Row 370: Fecha, Variedad, Calidad, Volumen, Precio minimo/maximo/promedio, Unidad, Precio $/Kg, Kg/unidad
$ws.Range("D370").Value = 44769
$ws.Range("K370").Value = 'Caramelo'
$ws.Range("L370").Value = 'Segunda'
$ws.Range("M370").Value = 100
$ws.Range("N370").Value = 19000
$ws.Range("O370").Value = 20000
$ws.Range("P370").Value = 19500
$ws.Range("Q370").Value = '$/caja 14 unidades'
$ws.Range("S370").Value = 1393
$ws.Range("T370").Value = 14

# Row 371: Fecha, Variedad, Calidad, Volumen, Precio minimo/maximo/promedio, Unidad, Precio $/Kg, Kg/unidad
$ws.Range("D371").Value = 44260
$ws.Range("K371").Value = 'Caramelo'
$ws.Range("L371").Value = 'Segunda'
$ws.Range("M371").Value = 100
$ws.Range("N371").Value = 15000
$ws.Range("O371").Value = 16000
$ws.Range("P371").Value = 15500
$ws.Range("Q371").Value = '$/caja 14 unidades'
$ws.Range("S371").Value = 1107
$ws.Range("T371").Value = 14
